# Actualiza base de datos EC: intercambia los valores de Mora entre la
# primera y la ultima fila de la tabla (F16 <-> F22).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F16").Value = 52000
$ws.Range("F22").Value = 24266
